$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Number of features: 31" -> "Number of features: 30"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Number of features: 31", $true, $false, $false, $false, $false, $true, 1, $false, "Number of features: 30", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Graduation rate of the class: 0.67%" -> "Graduation rate of the class: 67.09%"
#    The target splits the new number into separate runs ("67" / ".09") with
#    the _GoBack bookmark landing right after ".09" (before the trailing "%"),
#    mirroring how Word leaves runs split at the last edited/typed location.
# ---------------------------------------------------------------------------
$gradRange = $d.Content
$found = $gradRange.Find.Execute("0.67")
if ($found) {
    $startPos = $gradRange.Start

    # Replace "0.67" with the full new number; "%" just after stays untouched.
    $gradRange.Text = "67.09"

    $afterNumberPos = $startPos + 5      # right after "67.09", before "%"
    $midPos = $startPos + 2              # between "67" and ".09"

    # Move the document's _GoBack bookmark here (only one _GoBack may exist).
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $d.Range($afterNumberPos, $afterNumberPos))

    # Temporary bookmarks force the run splits we want; removing them again
    # afterwards leaves the runs split apart (no further text edit touches
    # them, so they are not re-coalesced).
    $d.Bookmarks.Add("zzzTempSplitA", $d.Range($midPos, $midPos))
    $d.Bookmarks.Add("zzzTempSplitB", $d.Range($startPos, $startPos))
    $d.Bookmarks("zzzTempSplitA").Delete()
    $d.Bookmarks("zzzTempSplitB").Delete()
}

# ---------------------------------------------------------------------------
# 3) Collapse the "Training Time (" / "secs" / ")" runs (with proofErr spell
#    markers around "secs") back into a single "Training Time (secs)" run.
#    Likewise for "Prediction Time (...)". Each occurs 3 times in the doc,
#    always alone in their own paragraph, so a whole-document replace-all
#    cannot bleed into neighboring runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Training Time (secs)", $true, $false, $false, $false, $false, $true, 1, $false, "Training Time (secs)", 2) | Out-Null
$d.Content.Find.Execute("Prediction Time (secs)", $true, $false, $false, $false, $false, $true, 1, $false, "Prediction Time (secs)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Merge the two runs that the old _GoBack bookmark used to split apart in
#    the closing paragraph ("...a reasonable" + " accuracy score...") back
#    into one run, and drop that bookmark (it has moved to the graduation
#    rate line above).
# ---------------------------------------------------------------------------
$tailText = "  If the test set produces a reasonable accuracy score we can begin the process of introducing new data and using the best logistic regression model to make real world predictions."
$tailRange = $d.Content
$foundTail = $tailRange.Find.Execute("  If the test set produces a reasonable")
if ($foundTail) {
    $tailStart = $tailRange.Start

    # Block the merge from bleeding into the preceding ")" run (which shares
    # identical run formatting) by marking the boundary temporarily.
    $d.Bookmarks.Add("zzzTempBlocker", $d.Range($tailStart, $tailStart))

    $mergeRange = $d.Range($tailStart, $tailStart + $tailText.Length)
    $mergeRange.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, $tailText, 2) | Out-Null

    $d.Bookmarks("zzzTempBlocker").Delete()
}
